$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3394
$ws.Range("F6").Value = 7916
$ws.Range("F9").Value = 12
$ws.Range("F12").Value = 176
$ws.Range("F13").Value = 1660
$ws.Range("F14").Value = 66
$ws.Range("F15").Value = 1101
$ws.Range("F16").Value = 1063
$ws.Range("F19").Value = 8557
$ws.Range("F20").Value = 205
$ws.Range("F25").Value = 1026
$ws.Range("F27").Value = 1187
$ws.Range("F32").Value = 113
$ws.Range("F33").Value = 1050
$ws.Range("F35").Value = 474
$ws.Range("F37").Value = 3554
$ws.Range("F39").Value = 47
$ws.Range("F41").Value = 767
$ws.Range("F43").Value = 118
$ws.Range("F45").Value = 678
$ws.Range("F48").Value = 25
$ws.Range("F49").Value = 2417

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 5
$ws.Range("F7").Value = 5
$ws.Range("F10").Value = 205
$ws.Range("F16").Value = 6
$ws.Range("C24").Value = "上海·幻彩唯响·森罗万象 Solo Live"
$ws.Range("F24").Value = 117
$ws.Range("F25").Value = 6956
$ws.Range("F33").Value = 38
$ws.Range("F37").Value = 3

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2081
$ws.Range("F5").Value = 1399
$ws.Range("F8").Value = 2248
$ws.Range("F9").Value = 9063
$ws.Range("F10").Value = 1333

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3394
$ws.Range("F4").Value = 2081
$ws.Range("F5").Value = 7916
$ws.Range("F6").Value = 2248
$ws.Range("F8").Value = 1333
$ws.Range("F11").Value = 12
$ws.Range("F13").Value = 1660
$ws.Range("F14").Value = 66
$ws.Range("F15").Value = 1101
$ws.Range("F16").Value = 1063
$ws.Range("F17").Value = 8557
$ws.Range("F18").Value = 205
$ws.Range("F22").Value = 1026
$ws.Range("F24").Value = 1187
$ws.Range("F27").Value = 205
$ws.Range("F29").Value = 113
$ws.Range("F30").Value = 1050
$ws.Range("F32").Value = 474
$ws.Range("F36").Value = 3554
$ws.Range("F37").Value = 47
$ws.Range("F38").Value = 767
$ws.Range("F41").Value = 678
$ws.Range("F46").Value = 25
$ws.Range("F48").Value = 2417
